$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.321.68'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.07%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.839.67'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.25%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.43%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6274'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.12%  '

# Row 7
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07428'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.74%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.95'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.57%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2890'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.37%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07722'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.11%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.832.67'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.66%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.952'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.01%  '

# Row 14
$ws.Range("E14").Value = '  -0.76%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001020'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.69%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.61'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.54%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.214'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.82%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.292.77'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.25%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.13'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.02%  '

# Row 20
$ws.Range("E20").Value = '  -0.50%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("D21").ClearFormats()

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.341'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.68%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.14%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.06'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.55%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.465'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.70%  '

# Row 26
$ws.Range("E26").Value = '  -2.36%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.34'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.19%  '

# Row 28
$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.07438'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +16.04%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.453'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.38%  '

# Row 30
$ws.Range("E30").Value = '  +0.40%  '

# Row 31
$ws.Range("E31").Value = '  -1.40%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.043'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.34%  '

# Row 33
$ws.Range("E33").Value = '  -0.40%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.138'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.12%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6936'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.60%  '

# Row 36
$ws.Range("E36").Value = '  -0.15%  '

# Row 37
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.949'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.48%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01837'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.56%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.813'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.56%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.233.25'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9336'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.49%  '

# Row 42
$ws.Range("E42").Value = '  +0.11%  '

# Row 43
$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.981.32'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.21%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.90'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.48%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.22'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.36%  '

# Row 46
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.703'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.46%  '

# Row 47
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.933'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.99%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1138'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.05%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.856'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.67%  '

# Row 50
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3896'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.03%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05665'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.69%  '
